# Scheduled data refresh: update market-board price/profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 197
$ws.Range("I4").Value = 197
$ws.Range("K4").Value = 197
$ws.Range("M4").Value = -83

$ws.Range("H107").Value = 1181.4286
$ws.Range("I107").Value = 1234.2307
$ws.Range("K107").Value = 1234.2307
$ws.Range("M107").Value = 685.7692999999999

$ws.Range("H137").Value = 1547.5641
$ws.Range("I137").Value = 1525.8846
$ws.Range("K137").Value = 4577.6538
$ws.Range("M137").Value = -2027.6538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1480.0448
$ws.Range("I32").Value = 1481.197
$ws.Range("K32").Value = 1481.197
$ws.Range("M32").Value = -1194.197

$ws.Range("H34").Value = 49998
$ws.Range("I34").Value = 49998
$ws.Range("K34").Value = 49998
$ws.Range("M34").Value = -49727

$ws.Range("H45").Value = 2378.25
$ws.Range("I45").Value = 2166.6667
$ws.Range("J45").Value = 3013
$ws.Range("K45").Value = 2166.6667
$ws.Range("L45").Value = 3013
$ws.Range("M45").Value = -1789.6667
$ws.Range("N45").Value = -3767

$ws.Range("H61").Value = 3364.8518
$ws.Range("I61").Value = 2734.8948
$ws.Range("K61").Value = 2734.8948
$ws.Range("M61").Value = -2522.8948

$ws.Range("H74").Value = 5788498.5
$ws.Range("I74").Value = 3088037
$ws.Range("K74").Value = 3088037
$ws.Range("M74").Value = -3087163

$ws.Range("H77").Value = 5788498.5
$ws.Range("I77").Value = 3088037
$ws.Range("K77").Value = 15440185
$ws.Range("M77").Value = -15435817

$ws.Range("H97").Value = 1537.0834
$ws.Range("J97").Value = 1980.5
$ws.Range("L97").Value = 1980.5
$ws.Range("N97").Value = -2972.5

$ws.Range("H110").Value = 1119.5714
$ws.Range("I110").Value = 821.0769
$ws.Range("K110").Value = 821.0769
$ws.Range("M110").Value = 1223.9231

$ws.Range("H136").Value = 3364.8518
$ws.Range("I136").Value = 2734.8948
$ws.Range("K136").Value = 8204.6844
$ws.Range("M136").Value = -5654.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1062.8334
$ws.Range("I5").Value = 1449.5
$ws.Range("K5").Value = 1449.5
$ws.Range("M5").Value = -1336.5

$ws.Range("H82").Value = 7160.125
$ws.Range("I82").Value = 4285.4287
$ws.Range("K82").Value = 4285.4287
$ws.Range("M82").Value = -3902.4287

$ws.Range("H85").Value = 7160.125
$ws.Range("I85").Value = 4285.4287
$ws.Range("K85").Value = 4285.4287
$ws.Range("M85").Value = -2959.4287

$ws.Range("H86").Value = 6511
$ws.Range("J86").Value = 9426.727999999999
$ws.Range("L86").Value = 9426.727999999999
$ws.Range("N86").Value = -11672.728

$ws.Range("H89").Value = 6511
$ws.Range("J89").Value = 9426.727999999999
$ws.Range("L89").Value = 47133.64
$ws.Range("N89").Value = -58365.64

$ws.Range("H94").Value = 1334.3334
$ws.Range("I94").Value = 1334.3334
$ws.Range("K94").Value = 1334.3334
$ws.Range("M94").Value = -883.3334

$ws.Range("H107").Value = 2408.7646
$ws.Range("J107").Value = 2929.5833
$ws.Range("L107").Value = 2929.5833
$ws.Range("N107").Value = -6769.5833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2486.5334
$ws.Range("I86").Value = 2330.2
$ws.Range("K86").Value = 2330.2
$ws.Range("M86").Value = -1207.2

$ws.Range("H89").Value = 2486.5334
$ws.Range("I89").Value = 2330.2
$ws.Range("K89").Value = 11651
$ws.Range("M89").Value = -6035

$ws.Range("H135").Value = 74654.336
$ws.Range("J135").Value = 74654.336
$ws.Range("L135").Value = 74654.336
$ws.Range("N135").Value = -84794.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 846303.4399999999
$ws.Range("I4").Value = 833495.4399999999
$ws.Range("K4").Value = 2500486.32
$ws.Range("M4").Value = -2500374.32

$ws.Range("H10").Value = 300.75
$ws.Range("I10").Value = 333.66666
$ws.Range("K10").Value = 1000.99998
$ws.Range("M10").Value = -861.9999799999999

$ws.Range("H11").Value = 511
$ws.Range("I11").Value = 462.5
$ws.Range("J11").Value = 899
$ws.Range("K11").Value = 1387.5
$ws.Range("L11").Value = 2697
$ws.Range("M11").Value = -1247.5
$ws.Range("N11").Value = -2977

$ws.Range("H18").Value = 2303
$ws.Range("I18").Value = 2163.8
$ws.Range("K18").Value = 6491.400000000001
$ws.Range("M18").Value = -6322.400000000001

$ws.Range("H46").Value = 1003927.56
$ws.Range("J46").Value = 5298.75
$ws.Range("L46").Value = 15896.25
$ws.Range("N46").Value = -16078.25

$ws.Range("H103").Value = 627.75
$ws.Range("I103").Value = 349
$ws.Range("J103").Value = 1018
$ws.Range("K103").Value = 1047
$ws.Range("L103").Value = 3054
$ws.Range("M103").Value = -168
$ws.Range("N103").Value = -4812

$ws.Range("H113").Value = 352.0435
$ws.Range("I113").Value = 289.6
$ws.Range("J113").Value = 369.3889
$ws.Range("K113").Value = 868.8000000000001
$ws.Range("L113").Value = 1108.1667
$ws.Range("M113").Value = 1301.2
$ws.Range("N113").Value = -5448.1667

$ws.Range("H130").Value = 15699.777
$ws.Range("I130").Value = 5023
$ws.Range("J130").Value = 17034.375
$ws.Range("K130").Value = 15069
$ws.Range("L130").Value = 51103.125
$ws.Range("N130").Value = -61143.125
$ws.Range("M130").Value = -10049

$ws.Range("H134").Value = 3214.2
$ws.Range("I134").Value = 3069.1177
$ws.Range("K134").Value = 9207.3531
$ws.Range("M134").Value = -4137.3531

$ws.Range("H136").Value = 6183.4
$ws.Range("I136").Value = 4994.5
$ws.Range("K136").Value = 14983.5
$ws.Range("M136").Value = -9883.5

$ws.Range("H138").Value = 8867.666999999999
$ws.Range("I138").Value = 6196.8
$ws.Range("K138").Value = 18590.4
$ws.Range("M138").Value = -13450.4

$ws.Range("H140").Value = 5948.6
$ws.Range("I140").Value = 6521.9414
$ws.Range("K140").Value = 19565.8242
$ws.Range("M140").Value = -14385.8242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2523.182
$ws.Range("I132").Value = 2470.25
$ws.Range("K132").Value = 7410.75
$ws.Range("M132").Value = -4880.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1470.4348
$ws.Range("I16").Value = 1391.15
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 1391.15
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -1221.15
$ws.Range("N16").Value = -2339

$ws.Range("H68").Value = 10392.934
$ws.Range("I68").Value = 11990.454
$ws.Range("J68").Value = 5999.75
$ws.Range("K68").Value = 11990.454
$ws.Range("L68").Value = 5999.75
$ws.Range("M68").Value = -11241.454
$ws.Range("N68").Value = -7497.75

$ws.Range("H71").Value = 10392.934
$ws.Range("I71").Value = 11990.454
$ws.Range("J71").Value = 5999.75
$ws.Range("K71").Value = 59952.27
$ws.Range("L71").Value = 29998.75
$ws.Range("M71").Value = -56208.27
$ws.Range("N71").Value = -37486.75

$ws.Range("H93").Value = 1511.8334
$ws.Range("I93").Value = 1663.1428
$ws.Range("J93").Value = 1300
$ws.Range("K93").Value = 1663.1428
$ws.Range("L93").Value = 1300
$ws.Range("M93").Value = -415.1428000000001
$ws.Range("N93").Value = -3796

$ws.Range("H99").Value = 30476
$ws.Range("J99").Value = 43979
$ws.Range("L99").Value = 43979
$ws.Range("N99").Value = -49969

$ws.Range("H136").Value = 43481280
$ws.Range("I136").Value = 3155.7273
$ws.Range("J136").Value = 1000000000
$ws.Range("K136").Value = 9467.1819
$ws.Range("L136").Value = 3000000000
$ws.Range("M136").Value = -6917.1819
$ws.Range("N136").Value = -3000005100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 38700
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H109").Value = 23226
$ws.Range("J109").Value = 23226
$ws.Range("L109").Value = 23226
$ws.Range("N109").Value = -26000

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = 0
